# Auto-generated edit script applying scheduled-runner Leve profit updates
# across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

function Set-CellValue {
    param($ws, $cellRef, $value)
    $ws.Range($cellRef).Value = $value
}

function Clear-CellValue {
    param($ws, $cellRef)
    $ws.Range($cellRef).ClearContents()
}

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
Set-CellValue $ws "H40" 3854.4443   # was 4224.875
Set-CellValue $ws "I40" 895.5   # was 900
Set-CellValue $ws "K40" 895.5   # was 900
Set-CellValue $ws "M40" -720.5   # was -725
Set-CellValue $ws "H55" 546.875   # was 449.9091
Set-CellValue $ws "I55" 460   # was 375
Set-CellValue $ws "J55" 599   # was 492.7143
Set-CellValue $ws "K55" 460   # was 375
Set-CellValue $ws "L55" 599   # was 492.7143
Set-CellValue $ws "M55" -246   # was -161
Set-CellValue $ws "N55" -1027   # was -920.7143
Set-CellValue $ws "H112" 2174.5789   # was 2146.15
Set-CellValue $ws "J112" 1850.9445   # was 1838.0526
Set-CellValue $ws "L112" 5552.833500000001   # was 5514.1578
Set-CellValue $ws "N112" -7768.833500000001   # was -7730.1578
Set-CellValue $ws "H137" 5861.488   # was 5977.923
Set-CellValue $ws "I137" 6008.905   # was 6263.421
Set-CellValue $ws "K137" 18026.715   # was 18790.263
Set-CellValue $ws "M137" -15476.715   # was -16240.263
Set-CellValue $ws "H141" 24937.13   # was 24977.105
Set-CellValue $ws "I141" 26295.584   # was 25638.918
Set-CellValue $ws "J141" 485   # was 490
Set-CellValue $ws "K141" 78886.75199999999   # was 76916.754
Set-CellValue $ws "L141" 1455   # was 1470
Set-CellValue $ws "M141" -73706.75199999999   # was -71736.754
Set-CellValue $ws "N141" -11815   # was -11830

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
Set-CellValue $ws "H32" 7647344   # was 7856857.5
Set-CellValue $ws "I32" 1250126   # was 1294769.8
Set-CellValue $ws "K32" 1250126   # was 1294769.8
Set-CellValue $ws "M32" -1249839   # was -1294482.8
Set-CellValue $ws "H37" 6262004   # was 6264625
Set-CellValue $ws "I37" 7152005   # was 8342500
Set-CellValue $ws "J37" 32000   # was 31000
Set-CellValue $ws "K37" 7152005   # was 8342500
Set-CellValue $ws "L37" 32000   # was 31000
Set-CellValue $ws "M37" -7151732   # was -8342227
Set-CellValue $ws "N37" -32546   # was -31546
Set-CellValue $ws "H61" 2041.3334   # was 2102.5881
Set-CellValue $ws "I61" 1988.5625   # was 2054.4666
Set-CellValue $ws "K61" 1988.5625   # was 2054.4666
Set-CellValue $ws "M61" -1776.5625   # was -1842.4666
Set-CellValue $ws "H74" 1958.7059   # was 2006.125
Set-CellValue $ws "I74" 2089.4614   # was 2163.5833
Set-CellValue $ws "K74" 2089.4614   # was 2163.5833
Set-CellValue $ws "M74" -1215.4614   # was -1289.5833
Set-CellValue $ws "H77" 1958.7059   # was 2006.125
Set-CellValue $ws "I77" 2089.4614   # was 2163.5833
Set-CellValue $ws "K77" 10447.307   # was 10817.9165
Set-CellValue $ws "M77" -6079.307000000001   # was -6449.916499999999
Set-CellValue $ws "H102" 1698.7693   # was 1613.9286
Set-CellValue $ws "J102" 1249.1666   # was 1143.7142
Set-CellValue $ws "L102" 1249.1666   # was 1143.7142
Set-CellValue $ws "N102" -4493.1666   # was -4387.7142
Set-CellValue $ws "H110" 1667.0454   # was 1704.5714
Set-CellValue $ws "I110" 1311.5294   # was 1338.5625
Set-CellValue $ws "K110" 1311.5294   # was 1338.5625
Set-CellValue $ws "M110" 733.4706000000001   # was 706.4375
Set-CellValue $ws "H132" 4496.8066   # was 4938.222
Set-CellValue $ws "I132" 4347.7036   # was 4839.9565
Set-CellValue $ws "K132" 13043.1108   # was 14519.8695
Set-CellValue $ws "M132" -10513.1108   # was -11989.8695
Set-CellValue $ws "H136" 2041.3334   # was 2102.5881
Set-CellValue $ws "I136" 1988.5625   # was 2054.4666
Set-CellValue $ws "K136" 5965.6875   # was 6163.399800000001
Set-CellValue $ws "M136" -3415.6875   # was -3613.399800000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
Set-CellValue $ws "H20" 3618.262   # was 3694.122
Set-CellValue $ws "I20" 3389.6296   # was 3500.4614
Set-CellValue $ws "K20" 3389.6296   # was 3500.4614
Set-CellValue $ws "M20" -3142.6296   # was -3253.4614
Set-CellValue $ws "H22" 0   # was 1000
Set-CellValue $ws "I22" 0   # was 1000
Set-CellValue $ws "K22" 0   # was 1000
Clear-CellValue $ws "M22"   # was -827
Set-CellValue $ws "H86" 981.4167   # was 1038.8182
Set-CellValue $ws "I86" 597.7   # was 625.2222
Set-CellValue $ws "K86" 597.7   # was 625.2222
Set-CellValue $ws "M86" 525.3   # was 497.7778
Set-CellValue $ws "H89" 981.4167   # was 1038.8182
Set-CellValue $ws "I89" 597.7   # was 625.2222
Set-CellValue $ws "K89" 2988.5   # was 3126.111
Set-CellValue $ws "M89" 2627.5   # was 2489.889
Set-CellValue $ws "H134" 3480.975   # was 3671.5405
Set-CellValue $ws "I134" 2634.3872   # was 2795.5
Set-CellValue $ws "K134" 7903.1616   # was 8386.5
Set-CellValue $ws "M134" -5368.1616   # was -5851.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
Set-CellValue $ws "H31" 3943.3215   # was 4234.52
Set-CellValue $ws "I31" 1346.5   # was 1489.375
Set-CellValue $ws "J31" 5386   # was 5526.353
Set-CellValue $ws "K31" 1346.5   # was 1489.375
Set-CellValue $ws "L31" 5386   # was 5526.353
Set-CellValue $ws "M31" -1051.5   # was -1194.375
Set-CellValue $ws "N31" -5976   # was -6116.353
Set-CellValue $ws "H34" 3943.3215   # was 4234.52
Set-CellValue $ws "I34" 1346.5   # was 1489.375
Set-CellValue $ws "J34" 5386   # was 5526.353
Set-CellValue $ws "K34" 1346.5   # was 1489.375
Set-CellValue $ws "L34" 5386   # was 5526.353
Set-CellValue $ws "M34" -1144.5   # was -1287.375
Set-CellValue $ws "N34" -5790   # was -5930.353
Set-CellValue $ws "H69" 5624.75   # was 5833.3335
Set-CellValue $ws "I69" 5624.75   # was 5833.3335
Set-CellValue $ws "K69" 5624.75   # was 5833.3335
Set-CellValue $ws "M69" -4875.75   # was -5084.3335
Set-CellValue $ws "H72" 5624.75   # was 5833.3335
Set-CellValue $ws "I72" 5624.75   # was 5833.3335
Set-CellValue $ws "K72" 16874.25   # was 17500.0005
Set-CellValue $ws "M72" -13130.25   # was -13756.0005
Set-CellValue $ws "H93" 12145.333   # was 13538.625
Set-CellValue $ws "I93" 6916   # was 7902.1665
Set-CellValue $ws "K93" 6916   # was 7902.1665
Set-CellValue $ws "M93" -5044   # was -6030.1665
Set-CellValue $ws "H103" 14098   # was 15918.4
Set-CellValue $ws "I103" 14864.333   # was 17523
Set-CellValue $ws "K103" 14864.333   # was 17523
Set-CellValue $ws "M103" -13692.333   # was -16351
Set-CellValue $ws "H132" 4017.647   # was 4081.75
Set-CellValue $ws "I132" 4056.5625   # was 4127.533
Set-CellValue $ws "K132" 12169.6875   # was 12382.599
Set-CellValue $ws "M132" -9639.6875   # was -9852.599000000002
Set-CellValue $ws "H134" 1798.075   # was 2026.1389
Set-CellValue $ws "I134" 1420.7646   # was 1609
Set-CellValue $ws "J134" 3936.1667   # was 4612.4
Set-CellValue $ws "K134" 4262.293799999999   # was 4827
Set-CellValue $ws "L134" 11808.5001   # was 13837.2
Set-CellValue $ws "M134" -1727.293799999999   # was -2292
Set-CellValue $ws "N134" -16878.5001   # was -18907.2
Set-CellValue $ws "H135" 0   # was 79999
Set-CellValue $ws "J135" 0   # was 79999
Set-CellValue $ws "L135" 0   # was 79999
Clear-CellValue $ws "N135"   # was -90139

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
Set-CellValue $ws "H10" 157.5   # was 102.5
Set-CellValue $ws "I10" 132.2   # was 100.111115
Set-CellValue $ws "J10" 199.66667   # was 104.888885
Set-CellValue $ws "K10" 396.6   # was 300.333345
Set-CellValue $ws "L10" 599.00001   # was 314.666655
Set-CellValue $ws "M10" -257.6   # was -161.333345
Set-CellValue $ws "N10" -877.00001   # was -592.666655
Set-CellValue $ws "H56" 8486.833000000001   # was 10258.692
Set-CellValue $ws "I56" 8486.833000000001   # was 10258.692
Set-CellValue $ws "K56" 8486.833000000001   # was 10258.692
Set-CellValue $ws "M56" -7956.833000000001   # was -9728.691999999999
Set-CellValue $ws "H137" 3574732.2   # was 3574792
Set-CellValue $ws "I137" 11112807   # was 12501742
Set-CellValue $ws "J137" 4065.4211   # was 4012.2
Set-CellValue $ws "K137" 33338421   # was 37505226
Set-CellValue $ws "L137" 12196.2633   # was 12036.6
Set-CellValue $ws "M137" -33333321   # was -37500126
Set-CellValue $ws "N137" -22396.2633   # was -22236.6

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
Set-CellValue $ws "H96" 30261   # was 0
Set-CellValue $ws "J96" 30261   # was 0
Set-CellValue $ws "L96" 30261   # was 0
Set-CellValue $ws "N96" -35753   # new cell (none previously)
Set-CellValue $ws "H107" 8821.058000000001   # was 9363.091
Set-CellValue $ws "I107" 12409.454   # was 13012.19
Set-CellValue $ws "J107" 2748.3845   # was 2977.1667
Set-CellValue $ws "K107" 12409.454   # was 13012.19
Set-CellValue $ws "L107" 2748.3845   # was 2977.1667
Set-CellValue $ws "M107" -10489.454   # was -11092.19
Set-CellValue $ws "N107" -6588.3845   # was -6817.1667
Set-CellValue $ws "H132" 3782.1091   # was 3932.151
Set-CellValue $ws "I132" 3674.5625   # was 3842.761
Set-CellValue $ws "K132" 11023.6875   # was 11528.283
Set-CellValue $ws "M132" -8493.6875   # was -8998.282999999999
Set-CellValue $ws "H136" 37540.5   # was 42714
Set-CellValue $ws "J136" 37540.5   # was 42714
Set-CellValue $ws "L136" 112621.5   # was 128142
Set-CellValue $ws "N136" -117721.5   # was -133242

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
Set-CellValue $ws "H22" 1060.4615   # was 1060.9259
Set-CellValue $ws "I22" 1033.3334   # was 1012.5
Set-CellValue $ws "J22" 1097.4546   # was 1131.3636
Set-CellValue $ws "K22" 1033.3334   # was 1012.5
Set-CellValue $ws "L22" 1097.4546   # was 1131.3636
Set-CellValue $ws "M22" -738.3334   # was -717.5
Set-CellValue $ws "N22" -1687.4546   # was -1721.3636
Set-CellValue $ws "H27" 1060.4615   # was 1060.9259
Set-CellValue $ws "I27" 1033.3334   # was 1012.5
Set-CellValue $ws "J27" 1097.4546   # was 1131.3636
Set-CellValue $ws "K27" 1033.3334   # was 1012.5
Set-CellValue $ws "L27" 1097.4546   # was 1131.3636
Set-CellValue $ws "M27" -926.3334   # was -905.5
Set-CellValue $ws "N27" -1311.4546   # was -1345.3636
Set-CellValue $ws "H100" 24006.646   # was 23559.693
Set-CellValue $ws "I100" 16745.805   # was 16397.238
Set-CellValue $ws "K100" 16745.805   # was 16397.238
Set-CellValue $ws "M100" -16204.805   # was -15856.238
Set-CellValue $ws "H132" 5438.25   # was 5800.6
Set-CellValue $ws "I132" 4802.2   # was 5002.6665
Set-CellValue $ws "J132" 6498.3335   # was 6997.5
Set-CellValue $ws "K132" 14406.6   # was 15007.9995
Set-CellValue $ws "L132" 19495.0005   # was 20992.5
Set-CellValue $ws "M132" -11876.6   # was -12477.9995
Set-CellValue $ws "N132" -24555.0005   # was -26052.5
Set-CellValue $ws "H136" 4366.923   # was 5163.3335
Set-CellValue $ws "I136" 3349.375   # was 4099
Set-CellValue $ws "J136" 5995   # was 6493.75
Set-CellValue $ws "K136" 10048.125   # was 12297
Set-CellValue $ws "L136" 17985   # was 19481.25
Set-CellValue $ws "M136" -7498.125   # was -9747
Set-CellValue $ws "N136" -23085   # was -24581.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
Set-CellValue $ws "H81" 2500550   # was 5000000
Set-CellValue $ws "I81" 2500550   # was 5000000
Set-CellValue $ws "K81" 5001100   # was 10000000
Set-CellValue $ws "M81" -5000039   # was -9998939
Set-CellValue $ws "H84" 2500550   # was 5000000
Set-CellValue $ws "I84" 2500550   # was 5000000
Set-CellValue $ws "K84" 25005500   # was 50000000
Set-CellValue $ws "M84" -25000196   # was -49994696
Set-CellValue $ws "H132" 11179.8125   # was 12690.929
Set-CellValue $ws "I132" 16962.8   # was 21053
Set-CellValue $ws "K132" 50888.39999999999   # was 63159
Set-CellValue $ws "M132" -48358.39999999999   # was -60629
Set-CellValue $ws "H136" 9432.058999999999   # was 10643.4
Set-CellValue $ws "I136" 6278.4165   # was 7464.7
Set-CellValue $ws "K136" 18835.2495   # was 22394.1
Set-CellValue $ws "M136" -16285.2495   # was -19844.1
